$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final values for A2:C21 after the edit (stock name lists reshuffled)
$data = @(
    @("立讯精密", "立讯精密", "山子高科"),
    @("和而泰", "山子高科", "立讯精密"),
    @("山子高科", "上海建工", "长盈精密"),
    @("东华软件", "和而泰", "和而泰"),
    @("凯美特气", "工业富联", "首开股份"),
    @("工业富联", "首开股份", "上海建工"),
    @("上海建工", "领益智造", "天普股份"),
    @("中科曙光", "东华软件", "工业富联"),
    @("长盈精密", "凯美特气", "领益智造"),
    @("首开股份", "中科曙光", "凯美特气"),
    @("福龙马", "长盈精密", "杭电股份"),
    @("领益智造", "福龙马", "东华软件"),
    @("卧龙电驱", "卧龙电驱", "中科曙光"),
    @("张江高科", "万马股份", "利欧股份"),
    @("初灵信息", "盈趣科技", "卧龙电驱"),
    @("万马股份", "歌尔股份", "万马股份"),
    @("英维克", "大众公用", "云南旅游"),
    @("利欧股份", "利欧股份", "山河智能"),
    @("三花智控", "张江高科", "福龙马"),
    @("杭电股份", "英维克", "三花智控")
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
